$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Excluded Studies" feature implementation:
# Replace the old "MM Maintenance" / "RRMM- Pfizer" category entries (rows 2 and 3,
# columns A/B) with the new "Test_Automation_1" / "Test_Automation_2" test categories.
$ws.Range("A2").Value = "Test_Automation_1"
$ws.Range("B2").Value = "Test_Automation_1_radio_button"
$ws.Range("A3").Value = "Test_Automation_2"
$ws.Range("B3").Value = "Test_Automation_2_radio_button"

# Move the active selection to B3, as reflected in the saved worksheet view.
[void]$ws.Range("B3").Select()
